# Apply update to the JenkinTrigger workbook:
#  - Reset AU_* rows (2-5) YES -> NO
#  - Set US_AccountCreate (row 18) and US_AccountUpdate (row 22) NO -> YES
#  - Append a new row 23: US_AccountDelete / US market DCS Account Delete / YES
#  - Extend data validation list and selection/view accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the AU rows from YES to NO
$ws.Range("C2").Value = "NO"
$ws.Range("C3").Value = "NO"
$ws.Range("C4").Value = "NO"
$ws.Range("C5").Value = "NO"

# Flip the US_AccountCreate / US_AccountUpdate rows from NO to YES
$ws.Range("C18").Value = "YES"
$ws.Range("C22").Value = "YES"

# Add the new US_AccountDelete row, reusing the formatting already used by
# the row above it (row 22) for columns B and C
$ws.Range("B22:C22").Copy()
$ws.Range("B23:C23").PasteSpecial(-4122)

$ws.Range("A23").Value = "US_AccountDelete"
$ws.Range("B23").Value = "US market DCS Account Delete"
$ws.Range("C23").Value = "YES"

# Give A23 a left/right border to match the rest of column A entries
$ws.Range("A23").Borders.Item(7).LineStyle = 1
$ws.Range("A23").Borders.Item(10).LineStyle = 1

# Extend the data validation range to include the new row
$ws.Range("C2:C23").Validation.Delete()
$ws.Range("C2:C23").Validation.Add(3, 1, 1, """YES,NO""")

# Update the view: scroll down a bit and select the newly-added cell
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A23").Select()
